$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,10
$row2[0,0] = -16.28558283343751
$row2[0,1] = -16.28558283343751
$row2[0,2] = -16.28558283343751
$row2[0,3] = -16.28558283343751
$row2[0,4] = -16.28558283343751
$row2[0,5] = -16.28558283343751
$row2[0,6] = -16.28558283343751
$row2[0,7] = -16.28558283343751
$row2[0,8] = -16.28558283343751
$row2[0,9] = -16.28558283343751
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object 'object[,]' 1,10
$row3[0,0] = -16.28558283343751
$row3[0,1] = -16.28558283343751
$row3[0,2] = -16.28558283343751
$row3[0,3] = -16.28558283343751
$row3[0,4] = -16.28558283343751
$row3[0,5] = -16.28558283343751
$row3[0,6] = -16.28558283343751
$row3[0,7] = 2.714618028315389
$row3[0,8] = -16.28558283343751
$row3[0,9] = -16.28558283343751
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object 'object[,]' 1,10
$row4[0,0] = -16.28558283343751
$row4[0,1] = -16.28558283343751
$row4[0,2] = 3.086876292590353
$row4[0,3] = -16.28558283343751
$row4[0,4] = 3.15674718255679
$row4[0,5] = -16.28558283343751
$row4[0,6] = 1.721208915999735
$row4[0,7] = -16.28558283343751
$row4[0,8] = 2.171657751856103
$row4[0,9] = -16.28558283343751
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object 'object[,]' 1,10
$row5[0,0] = -16.28558283343751
$row5[0,1] = -16.28558283343751
$row5[0,2] = -16.28558283343751
$row5[0,3] = -16.28558283343751
$row5[0,4] = -16.28558283343751
$row5[0,5] = 2.915164940338256
$row5[0,6] = -16.28558283343751
$row5[0,7] = -16.28558283343751
$row5[0,8] = -16.28558283343751
$row5[0,9] = -16.28558283343751
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object 'object[,]' 1,10
$row6[0,0] = -16.28558283343751
$row6[0,1] = -16.28558283343751
$row6[0,2] = -16.28558283343751
$row6[0,3] = -16.28558283343751
$row6[0,4] = -16.28558283343751
$row6[0,5] = -16.28558283343751
$row6[0,6] = -16.28558283343751
$row6[0,7] = -16.28558283343751
$row6[0,8] = -16.28558283343751
$row6[0,9] = -16.28558283343751
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object 'object[,]' 1,10
$row7[0,0] = 2.610941745126737
$row7[0,1] = -16.28558283343751
$row7[0,2] = -16.28558283343751
$row7[0,3] = -16.28558283343751
$row7[0,4] = -16.28558283343751
$row7[0,5] = -16.28558283343751
$row7[0,6] = -16.28558283343751
$row7[0,7] = -16.28558283343751
$row7[0,8] = -16.28558283343751
$row7[0,9] = -16.28558283343751
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object 'object[,]' 1,10
$row8[0,0] = -16.28558283343751
$row8[0,1] = -16.28558283343751
$row8[0,2] = -16.28558283343751
$row8[0,3] = 1.589140421182108
$row8[0,4] = -16.28558283343751
$row8[0,5] = -16.28558283343751
$row8[0,6] = -16.28558283343751
$row8[0,7] = -16.28558283343751
$row8[0,8] = -16.28558283343751
$row8[0,9] = -16.28558283343751
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object 'object[,]' 1,10
$row9[0,0] = 3.796052732472259
$row9[0,1] = -16.28558283343751
$row9[0,2] = -16.28558283343751
$row9[0,3] = -16.28558283343751
$row9[0,4] = -16.28558283343751
$row9[0,5] = -16.28558283343751
$row9[0,6] = -16.28558283343751
$row9[0,7] = -16.28558283343751
$row9[0,8] = -16.28558283343751
$row9[0,9] = -16.28558283343751
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object 'object[,]' 1,10
$row10[0,0] = -16.28558283343751
$row10[0,1] = -16.28558283343751
$row10[0,2] = -16.28558283343751
$row10[0,3] = -16.28558283343751
$row10[0,4] = -16.28558283343751
$row10[0,5] = -16.28558283343751
$row10[0,6] = -16.28558283343751
$row10[0,7] = 1.248241158711078
$row10[0,8] = -16.28558283343751
$row10[0,9] = 2.011212415078234
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object 'object[,]' 1,10
$row11[0,0] = -16.28558283343751
$row11[0,1] = -16.28558283343751
$row11[0,2] = -16.28558283343751
$row11[0,3] = 3.098013890691628
$row11[0,4] = -16.28558283343751
$row11[0,5] = 2.760008871985114
$row11[0,6] = -16.28558283343751
$row11[0,7] = -16.28558283343751
$row11[0,8] = -16.28558283343751
$row11[0,9] = 1.890536735761882
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object 'object[,]' 1,10
$row12[0,0] = -16.28558283343751
$row12[0,1] = -16.28558283343751
$row12[0,2] = -16.28558283343751
$row12[0,3] = -16.28558283343751
$row12[0,4] = -16.28558283343751
$row12[0,5] = -16.28558283343751
$row12[0,6] = -16.28558283343751
$row12[0,7] = -16.28558283343751
$row12[0,8] = -16.28558283343751
$row12[0,9] = -16.28558283343751
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object 'object[,]' 1,10
$row13[0,0] = -16.28558283343751
$row13[0,1] = -16.28558283343751
$row13[0,2] = -16.28558283343751
$row13[0,3] = 2.327375627962638
$row13[0,4] = -16.28558283343751
$row13[0,5] = -16.28558283343751
$row13[0,6] = -16.28558283343751
$row13[0,7] = -16.28558283343751
$row13[0,8] = 1.923189729324837
$row13[0,9] = 1.957446999803929
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object 'object[,]' 1,10
$row14[0,0] = -16.28558283343751
$row14[0,1] = -16.28558283343751
$row14[0,2] = 1.368528702041514
$row14[0,3] = -16.28558283343751
$row14[0,4] = -16.28558283343751
$row14[0,5] = -16.28558283343751
$row14[0,6] = -16.28558283343751
$row14[0,7] = -16.28558283343751
$row14[0,8] = -16.28558283343751
$row14[0,9] = 2.069944631311495
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object 'object[,]' 1,10
$row15[0,0] = -16.28558283343751
$row15[0,1] = -16.28558283343751
$row15[0,2] = 0.9774020405425654
$row15[0,3] = -16.28558283343751
$row15[0,4] = -16.28558283343751
$row15[0,5] = -16.28558283343751
$row15[0,6] = -16.28558283343751
$row15[0,7] = -16.28558283343751
$row15[0,8] = -16.28558283343751
$row15[0,9] = -16.28558283343751
$ws.Range("B15:K15").Value = $row15

$row16 = New-Object 'object[,]' 1,10
$row16[0,0] = -16.28558283343751
$row16[0,1] = -16.28558283343751
$row16[0,2] = -16.28558283343751
$row16[0,3] = -16.28558283343751
$row16[0,4] = -16.28558283343751
$row16[0,5] = -16.28558283343751
$row16[0,6] = -16.28558283343751
$row16[0,7] = -16.28558283343751
$row16[0,8] = 2.197244814026567
$row16[0,9] = -16.28558283343751
$ws.Range("B16:K16").Value = $row16

$row17 = New-Object 'object[,]' 1,10
$row17[0,0] = -16.28558283343751
$row17[0,1] = -16.28558283343751
$row17[0,2] = 1.677725330963569
$row17[0,3] = -16.28558283343751
$row17[0,4] = -16.28558283343751
$row17[0,5] = -16.28558283343751
$row17[0,6] = 1.495142743896241
$row17[0,7] = 1.944367215111849
$row17[0,8] = 2.038725188280642
$row17[0,9] = -16.28558283343751
$ws.Range("B17:K17").Value = $row17

$row18 = New-Object 'object[,]' 1,10
$row18[0,0] = -16.28558283343751
$row18[0,1] = -16.28558283343751
$row18[0,2] = -16.28558283343751
$row18[0,3] = -16.28558283343751
$row18[0,4] = -16.28558283343751
$row18[0,5] = -16.28558283343751
$row18[0,6] = 1.510057352014554
$row18[0,7] = 1.201260042497898
$row18[0,8] = 1.588299798254502
$row18[0,9] = -16.28558283343751
$ws.Range("B18:K18").Value = $row18

$row19 = New-Object 'object[,]' 1,10
$row19[0,0] = -16.28558283343751
$row19[0,1] = -16.28558283343751
$row19[0,2] = 1.098399833985733
$row19[0,3] = -16.28558283343751
$row19[0,4] = -16.28558283343751
$row19[0,5] = -16.28558283343751
$row19[0,6] = 1.610826443406027
$row19[0,7] = 1.4572331519074
$row19[0,8] = -16.28558283343751
$row19[0,9] = -16.28558283343751
$ws.Range("B19:K19").Value = $row19

$row20 = New-Object 'object[,]' 1,10
$row20[0,0] = -16.28558283343751
$row20[0,1] = 4.321910937474696
$row20[0,2] = 0.688753327412148
$row20[0,3] = -16.28558283343751
$row20[0,4] = 3.470095317439287
$row20[0,5] = -16.28558283343751
$row20[0,6] = 1.924617617285345
$row20[0,7] = 1.115171548737728
$row20[0,8] = -16.28558283343751
$row20[0,9] = 2.062959474711402
$ws.Range("B20:K20").Value = $row20

$row21 = New-Object 'object[,]' 1,10
$row21[0,0] = -16.28558283343751
$row21[0,1] = -16.28558283343751
$row21[0,2] = -16.28558283343751
$row21[0,3] = 1.769683712330965
$row21[0,4] = -16.28558283343751
$row21[0,5] = 2.506562439531462
$row21[0,6] = 2.065205311518139
$row21[0,7] = -16.28558283343751
$row21[0,8] = -16.28558283343751
$row21[0,9] = -16.28558283343751
$ws.Range("B21:K21").Value = $row21

